$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.875.33'
$ws.Range("D3").Value = '1.705.62'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.50'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3953'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4058'
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.484'
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.000'
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.44'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08814'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.42'
$ws.Range("E13").Value = '  +8.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.494'
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.127'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001360'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '1.698.92'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.70'
$ws.Range("E18").Value = '  -3.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07155'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +4.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.292'
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.37'
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("D24").Value = '24.872.16'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.988'
$ws.Range("E25").Value = '  -3.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.335'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.20'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.176'
$ws.Range("E28").Value = '  +18.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.70'
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.494'
$ws.Range("E30").Value = '  -8.99%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '145.10'
$ws.Range("E31").Value = '  +3.85%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.263'
$ws.Range("E32").Value = '  +15.07%  '
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.886.82'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08790'
$ws.Range("E34").Value = '  -4.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03205'
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.180'
$ws.Range("E36").Value = '  -10.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.035'
$ws.Range("E37").Value = '  -4.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2884'
$ws.Range("E38").Value = '  +2.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.89'
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8415'
$ws.Range("E40").Value = '  +7.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09260'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.472'
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("E44").Value = '  +7.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.689'
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7422'
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.390'
$ws.Range("E48").Value = '  +2.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9988'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.00'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08348'
$ws.Range("E51").Value = '  +3.42%  '
